$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 so the new header cells (I1, J1) match
# the existing bold/bordered/centered header formatting.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New column data (I = I0, J = IF)
$values = @{
    2  = @(7, 7)
    3  = @(5, 6)
    4  = @(6, 7)
    5  = @(8, 9)
    6  = @(8, 8)
    7  = @(9, 9)
    8  = @(8, 9)
    9  = @(6, 6)
    10 = @(6, 6)
    11 = @(7, 7)
    12 = @(5, 5)
    13 = @(5, 5)
    14 = @(8, 8)
    15 = @(6, 6)
    16 = @(9, 9)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
